$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing demand values ---
# Row 11: SE04 / National Trends->Distributed Energy / 2040   500 -> 800
$ws.Cells.Item(11, 5).Value = 800

# Row 14: DKW1 / Distributed Energy / 2040   500 -> 1000
$ws.Cells.Item(14, 5).Value = 1000

# Helper: give a cell the same "data column" look used by the rest of the
# Scenario/Year columns (bordered, bold Calibri, centered) so the new rows
# read consistently with the existing table.
function Set-DataColumnStyle($cell) {
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# --- Append new FI00 / industry rows (29-31) ---
$ws.Cells.Item(29, 1).Value = "FI00"
$ws.Cells.Item(29, 2).Value = "industry"
$ws.Cells.Item(29, 3).Value = "National Trends"
$ws.Cells.Item(29, 4).Value = 2025
$ws.Cells.Item(29, 5).Value = 1450
Set-DataColumnStyle($ws.Cells.Item(29, 3))
Set-DataColumnStyle($ws.Cells.Item(29, 4))

$ws.Cells.Item(30, 1).Value = "FI00"
$ws.Cells.Item(30, 2).Value = "industry"
$ws.Cells.Item(30, 3).Value = "Distributed Energy"
$ws.Cells.Item(30, 4).Value = 2030
$ws.Cells.Item(30, 5).Value = 1450
Set-DataColumnStyle($ws.Cells.Item(30, 3))

$ws.Cells.Item(31, 1).Value = "FI00"
$ws.Cells.Item(31, 2).Value = "industry"
$ws.Cells.Item(31, 3).Value = "Distributed Energy"
$ws.Cells.Item(31, 4).Value = 2040
$ws.Cells.Item(31, 5).Value = 1450
Set-DataColumnStyle($ws.Cells.Item(31, 3))
Set-DataColumnStyle($ws.Cells.Item(31, 4))

# --- Move active selection, matching the saved cursor position ---
$ws.Range("F30").Select()
